$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Checked the architecture document / user stories document / installation
#     manual columns (D = Architekturdokument, E = Installationsanleitung,
#     F = User Stories Dokument) ---

# User Stories Dokument (column F) results for several checks
$ws.Range("F8").Value  = "nicht nötig"
$ws.Range("F9").Value  = "nicht nötig"
$ws.Range("F11").Value = "keine schachtelung nötig"
$ws.Range("F14").Value = "weder noch vorhanden"
$ws.Range("F15").Value = "nicht vorhanden"
$ws.Range("F16").Value = "nicht vorhanden"
$ws.Range("F17").Value = "nicht vorhanden"
$ws.Range("F18").Value = "gecheckt am 17.06.2016"

# Installationsanleitung (column E) corrections
$ws.Range("E16").Value = "ja"
$ws.Range("E18").Value = "gecheckt am 17.06.16"

# New "Dokumentenstatus" row for all three documents
$ws.Range("D19").Value = "Dokumentenstatus: to be reviewed"
$ws.Range("E19").Value = "Dokumentenstatus: final"
$ws.Range("F19").Value = "Dokumentenstatus: to be reviewed"

# --- Un-hide + resize the previously-hidden helper columns C:E and refresh
#     the best-fit widths now that longer text has been entered ---
$ws.Columns.Item(4).Hidden = $false
$ws.Columns.Item(5).Hidden = $false
$ws.Columns.Item(3).ColumnWidth = 17.17

$ws.Range("D1:F20").EntireColumn.AutoFit() | Out-Null

# --- Selection left where the user ended up after the edit ---
$ws.Range("G6").Select()
